$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 31 (shifts existing rows 31-34 down to 32-35)
$ws.Rows("31:31").Insert()

# Populate the new row 31 with the new street "COMPLEXO VIARIO MARIA MALUF"
$ws.Cells.Item(31, 1).Value = "COMPLEXO VIARIO MARIA MALUF"
$ws.Cells.Item(31, 2).Value = "COMPLEXO VIARIO MARIA MALUF"
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(31, 4).Value = 6
$ws.Cells.Item(31, 5).Value = 2024
$ws.Cells.Item(31, 9).Value = 26
$ws.Cells.Item(31, 10).Value = 2

# Renumber id_logradouro (column I) for the rows that were pushed down
$ws.Cells.Item(32, 9).Value = 27
$ws.Cells.Item(33, 9).Value = 28
$ws.Cells.Item(34, 9).Value = 29
$ws.Cells.Item(35, 9).Value = 30

# Update the frozen-pane top-left cell and the active selection to match
$win = $excel.ActiveWindow
$win.FreezePanes = $false
[void]$ws.Range("A2").Select()
$win.FreezePanes = $true
[void]$ws.Range("J32").Select()
